# The workbook holds 16 "backward elimination" OLS summary sheets (tabs
# "18" down to "3"), each with the statsmodels text-report pasted into
# cell B2. The report was regenerated, so every sheet's timestamp moved
# from Sun, 22 Dec 2019 23:07:16 to Wed, 25 Dec 2019 23:10:04 (the
# regression numbers themselves are unchanged).

$wb = $excel.ActiveWorkbook

$oldDate = "Sun, 22 Dec 2019"
$newDate = "Wed, 25 Dec 2019"
$oldTime = "23:07:16"
$newTime = "23:10:04"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Cells.Item(2, 2)
    $text = $cell.Value2

    if ($text -ne $null -and $text -ne "") {
        $updated = $text.Replace($oldDate, $newDate).Replace($oldTime, $newTime)
        if ($updated -ne $text) {
            $cell.Value2 = $updated
        }
    }
}
